$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    ,@("Sunday, Jan 15", "2:10 AM", "E47830", "Marsa Alam", "(RMF)", "Enter Air ", "B738", "(SP-ENL)", "2:06 AM", $null, "0 hours, -4 minutes")
    ,@("Sunday, Jan 15", "10:55 AM", "FR5074", "Birmingham", "(BHX)", "Ryanair ", "B738", "(SP-RKR)", "10:45 AM", $null, "0 hours, -10 minutes")
    ,@("Sunday, Jan 15", "12:20 PM", "LH1638", "Munich", "(MUC)", "Lufthansa ", "CRJ9", "(D-ACNH)", "12:04 PM", $null, "0 hours, -16 minutes")
    ,@("Sunday, Jan 15", "1:30 PM", "FR7907", "Alicante", "(ALC)", "Ryanair ", "B738", "(SP-RSM)", "1:06 PM", $null, "0 hours, -24 minutes")
    ,@("Sunday, Jan 15", "2:15 PM", "LO3943", "Warsaw", "(WAW)", "LOT ", "E170", "(SP-LDI)", "2:13 PM", $null, "0 hours, -2 minutes")
    ,@("Sunday, Jan 15", "2:35 PM", "KL1273", "Amsterdam", "(AMS)", "KLM ", "E75L", "(PH-EXR)", "2:08 PM", $null, "0 hours, -27 minutes")
    ,@("Sunday, Jan 15", "2:50 PM", "LH1390", "Frankfurt", "(FRA)", "Lufthansa ", "CRJ9", "(D-ACNW)", "2:38 PM", $null, "0 hours, -12 minutes")
    ,@("Sunday, Jan 15", "4:00 PM", "W91902", "London", "(LTN)", "Wizz Air ", "A321", "(G-WUKI)", "3:36 PM", $null, "0 hours, -24 minutes")
    ,@("Sunday, Jan 15", "4:25 PM", "FR1750", "London", "(STN)", "Ryanair ", "B738", "(SP-RKR)", "4:12 PM", $null, "0 hours, -13 minutes")
)

$startRow = 159
$num = 158
foreach ($row in $newRows) {
    $ws.Cells.Item($startRow, 1).Value = $num
    for ($i = 0; $i -lt $row.Length; $i++) {
        $col = 2 + $i
        $val = $row[$i]
        if ($null -ne $val) {
            $ws.Cells.Item($startRow, $col).Value = $val
        }
    }
    $startRow++
    $num++
}
